# "Update countries & provincias Spain"
#
# Refreshes the COVID "Pais" snapshot (sheet "Pais"):
#   - bumps the "Datos actualizados ..." timestamp in A1 from 15:30 to 16:47
#   - refreshes totals (Casos totales/Nuevos casos/Casos activos/Recuperados/
#     Casos criticos/Muertes hoy/Muertes) for several countries
#   - "Reunion" moves up in the (case-count sorted) table to sit right after
#     "Islandia" and ahead of "Sierra Leona", which pushes Sierra Leona,
#     Trinidad yTobago, Malta and Botsuana each down one row (each of those
#     four keeps its own previous totals - only Islandia and Reunion get
#     genuinely new numbers)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 16:47"

# Estados Unidos
$ws.Range("B4").Value = 6392303
$ws.Range("C4").Value = 3246
$ws.Range("D4").Value = 3636310
$ws.Range("E4").Value = 2563825
$ws.Range("G4").Value = 57
$ws.Range("H4").Value = 192168

# India
$ws.Range("B6").Value = 4047653
$ws.Range("C6").Value = 27414
$ws.Range("D6").Value = 3124787
$ws.Range("E6").Value = 852946
$ws.Range("G6").Value = 285
$ws.Range("H6").Value = 69920

# Irak
$ws.Range("B23").Value = 256719
$ws.Range("C23").Value = 4644
$ws.Range("D23").Value = 195259
$ws.Range("E23").Value = 54038
$ws.Range("G23").Value = 63
$ws.Range("H23").Value = 7422

# Portugal
$ws.Range("B51").Value = 59943
$ws.Range("C51").Value = 486
$ws.Range("D51").Value = 42793
$ws.Range("E51").Value = 15312
$ws.Range("G51").Value = 5
$ws.Range("H51").Value = 1838

# Kenia
$ws.Range("B68").Value = 35020
$ws.Range("C68").Value = 136
$ws.Range("D68").Value = 21158
$ws.Range("E68").Value = 13268
$ws.Range("G68").Value = 5
$ws.Range("H68").Value = 594

# Republica de Macedonia
$ws.Range("B86").Value = 14998
$ws.Range("C86").Value = 127
$ws.Range("D86").Value = 12149
$ws.Range("E86").Value = 2235
$ws.Range("G86").Value = 5
$ws.Range("H86").Value = 614

# Zambia
$ws.Range("B89").Value = 12709
$ws.Range("C89").Value = 70
$ws.Range("D89").Value = 11668
$ws.Range("E89").Value = 749

# Noruega
$ws.Range("B91").Value = 11254
$ws.Range("C91").Value = 23
$ws.Range("E91").Value = 1642

# Islandia (row 143) - new totals, stays in place
$ws.Range("B143").Value = 2136
$ws.Range("C143").Value = 1
$ws.Range("D143").Value = 2038
$ws.Range("E143").Value = 88

# Reunion moves into row 144 (right after Islandia) with its new totals
$ws.Range("A144").Value = "Reunion"
$ws.Range("B144").Value = 2115
$ws.Range("C144").Value = 113
$ws.Range("D144").Value = 1313
$ws.Range("E144").Value = 791
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 11

# Sierra Leona shifts down to row 145, keeping its previous totals
$ws.Range("A145").Value = "Sierra Leona"
$ws.Range("B145").Value = 2041
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 1602
$ws.Range("E145").Value = 368
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 71

# Trinidad yTobago shifts down to row 146, keeping its previous totals
$ws.Range("A146").Value = "Trinidad yTobago"
$ws.Range("B146").Value = 2040
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 707
$ws.Range("E146").Value = 1302
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 31

# Malta shifts down to row 147, keeping its previous totals
$ws.Range("A147").Value = "Malta"
$ws.Range("B147").Value = 2014
$ws.Range("C147").Value = 30
$ws.Range("D147").Value = 1601
$ws.Range("E147").Value = 399
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 14

# Botsuana shifts down to row 148, keeping its previous totals
$ws.Range("A148").Value = "Botsuana"
$ws.Range("B148").Value = 2002
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 493
$ws.Range("E148").Value = 1501
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 8

# Birmania
$ws.Range("B160").Value = 1253
$ws.Range("C160").Value = 82
$ws.Range("D160").Value = 371
$ws.Range("E160").Value = 875
